$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the StageMaster / ItemMaster blocks down by 6 rows to make room for
# a new "ObjectMaster" master-data block (header + column-header + 2 data
# rows + a blank separator) inserted right after the PlayerMaster block.
$ws.Rows("8:13").Insert()

# Inserting rows duplicates the formerly-last blank spacer row (old row 17,
# "customFormat" only) down to row 23 as a now-redundant trailing row -
# drop it so the sheet ends cleanly at row 22, same as before the edit.
$ws.Rows("23:23").Delete()

# --- New "ObjectMaster" block (rows 8-11) -------------------------------
$ws.Range("A8").Value = "MASTERDATA"
$ws.Range("B8").Value = "ObjectMaster"

$ws.Range("A9").Value = "Code"
$ws.Range("B9").Value = "ObjectType"
$ws.Range("C9").Value = "ModelPrefabPath"
$ws.Range("D9").Value = "InitialPosition.x"
$ws.Range("E9").Value = "InitialPosition.y"
$ws.Range("F9").Value = "InitialPosition.z"

$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Rock"
$ws.Range("C10").Value = "SimpleNaturePack/Prefabs/Rock_05"
$ws.Range("D10").Value = -1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0

$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "Mushroom"
$ws.Range("C11").Value = "SimpleNaturePack/Prefabs/Mushroom_02"
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = -1

# --- Restore the saved view/selection state -----------------------------
$ws.Range("F12").Select()
